# UserStories.xlsx -- "Booking.cs und IBookingRepository.cs erstellt"
#
# On the "Tasks" sheet:
#  - fix the typo "Interface IBookingReposittory erstellen"
#    -> "Interface IBookingRepository erstellen" (row 41, column B)
#  - mark the three now-finished tasks (rows 39-41) as "done" with a
#    completion date of 26.03.2019 (Excel serial 43550), copying the
#    date formatting already used elsewhere in column D
#  - update the view: scrolled position / active selection moved from
#    D43 to E43

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# --- fix the spelling of the interface name -------------------------
$ws.Range("B41").Value = "Interface IBookingRepository erstellen"

# --- mark "Tabelle bookings und rooms erstellen" / "Klasse Booking
#     erstellen" / "Interface IBookingRepository erstellen" as done,
#     each with a completion date -------------------------------------
$ws.Range("C39").Value = "done"
$ws.Range("C40").Value = "done"
$ws.Range("C41").Value = "done"

# Copy the date-cell formatting already used in the sheet (e.g. D4)
# onto D39:D41 before writing the serial date values, so the new cells
# pick up the same number format / wrap-text style instead of Excel's
# generic default.
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D39:D41").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("D39").Value = 43550
$ws.Range("D40").Value = 43550
$ws.Range("D41").Value = 43550

# --- update the sheet's active selection/scroll position ------------
$ws.Select() | Out-Null
$ws.Range("E43").Select() | Out-Null
